$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 11 for year 2021
$ws.Range("A11").Value = "2021年"
$ws.Range("A10").Copy()
$ws.Range("A11").PasteSpecial(-4122)  # xlPasteFormats, to match bold/centered/border style used for year cells

$ws.Range("B11").Value = 2084.99
$ws.Range("C11").Value = 546.7
$ws.Range("D11").Value = 77.58
$ws.Range("F11").Value = 709.65
$ws.Range("G11").Value = 2480.22
$ws.Range("H11").Value = 132.41
$ws.Range("I11").Value = 1314.64
$ws.Range("J11").Value = 253.32
$ws.Range("K11").Value = 51032.69
$ws.Range("L11").Value = 246.87
$ws.Range("M11").Value = 30.68
$ws.Range("N11").Value = 27.91
$ws.Range("O11").Value = 467.41
$ws.Range("P11").Value = 615.55
$ws.Range("Q11").Value = 3.11
$ws.Range("R11").Value = 88.76000000000001
$ws.Range("S11").Value = 1278.07
$ws.Range("T11").Value = 148.71
$ws.Range("U11").Value = 6516.48
$ws.Range("W11").Value = 97.06999999999999
$ws.Range("X11").Value = 395.8
$ws.Range("Y11").Value = 1061.76
$ws.Range("Z11").Value = 4133.07
$ws.Range("AA11").Value = 344.48
$ws.Range("AB11").Value = 196.19
$ws.Range("AC11").Value = 179.75
$ws.Range("AD11").Value = 685.46
$ws.Range("AE11").Value = 598.4400000000001
$ws.Range("AF11").Value = 17560.93
$ws.Range("AG11").Value = 3142.53
$ws.Range("AH11").Value = 781.34
$ws.Range("AI11").Value = 384.18
$ws.Range("AJ11").Value = 114.92
$ws.Range("AK11").Value = 1477.91
$ws.Range("AL11").Value = 564.86
$ws.Range("AM11").Value = 1006.96
$ws.Range("AN11").Value = 12.82
$ws.Range("AO11").Value = 859.6799999999999
$ws.Range("AP11").Value = 368.63
$ws.Range("AQ11").Value = 42.13
